$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows from bottom to top to avoid index shifting issues.
# Old rows 80 (UCS 921) and 82 (UCS 941) are removed outright.
$ws.Rows("82:82").Delete()
$ws.Rows("80:80").Delete()
# Old rows 74-76 (UCS 621, UCS 622, UCS 623) are removed outright.
$ws.Rows("74:76").Delete()

# The remaining rows (old 77,78,79,81,83) have shifted up to become
# rows 74-78. Update their Model-name cells (column A) to the new
# "UCS-4T ..." labels and give them the header-style formatting (style
# of A1), matching the target workbook.
$ws.Range("A1").Copy()
$ws.Range("A74:A78").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A74").Value = "UCS-4T 641"
$ws.Range("A75").Value = "UCS-4T 642"
$ws.Range("A76").Value = "UCS-4T 643"
$ws.Range("A77").Value = "UCS-4T 941"
$ws.Range("A78").Value = "UCS-4T 942"

# Match the final on-screen selection/scroll position recorded in the
# saved workbook (cell F78, scrolled so row 61 is at the top).
[void]$ws.Range("F78").Select()
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "done"
